$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 17
$ws.Range("H17").Value = 1391473.2
$ws.Range("J17").Value = 1391473.2
$ws.Range("L17").Value = 4174419.6
$ws.Range("N17").Value = -4174755.6

# row 129
$ws.Range("H129").Value = 1249.9584
$ws.Range("I129").Value = 399.9375
$ws.Range("J129").Value = 2950
$ws.Range("K129").Value = 1199.8125
$ws.Range("L129").Value = 8850
$ws.Range("M129").Value = 3800.1875
$ws.Range("N129").Value = -18850

# row 137
$ws.Range("H137").Value = 5000774
$ws.Range("I137").Value = 692.44446
$ws.Range("K137").Value = 2077.33338
$ws.Range("M137").Value = 472.66662

# row 138
$ws.Range("H138").Value = 1457.9552
$ws.Range("I138").Value = 1176.5962
$ws.Range("J138").Value = 2433.3333
$ws.Range("K138").Value = 3529.7886
$ws.Range("L138").Value = 7299.999899999999
$ws.Range("M138").Value = 1610.2114
$ws.Range("N138").Value = -17579.9999

# row 139
$ws.Range("H139").Value = 162245
$ws.Range("J139").Value = 162245
$ws.Range("L139").Value = 162245
$ws.Range("N139").Value = -172525

# row 141
$ws.Range("H141").Value = 884.2632
$ws.Range("I141").Value = 856.14545
$ws.Range("J141").Value = 1657.5
$ws.Range("K141").Value = 2568.43635
$ws.Range("L141").Value = 4972.5
$ws.Range("M141").Value = 2611.56365
$ws.Range("N141").Value = -15332.5

$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 8544.188
$ws.Range("I32").Value = 9078.4
$ws.Range("J32").Value = 6808
$ws.Range("K32").Value = 9078.4
$ws.Range("L32").Value = 6808
$ws.Range("M32").Value = -8791.4
$ws.Range("N32").Value = -7382

# row 61
$ws.Range("H61").Value = 13159577
$ws.Range("I61").Value = 15153317
$ws.Range("J61").Value = 892.4
$ws.Range("K61").Value = 15153317
$ws.Range("L61").Value = 892.4
$ws.Range("M61").Value = -15153105
$ws.Range("N61").Value = -1316.4

# row 74
$ws.Range("H74").Value = 9092324
$ws.Range("I74").Value = 11906162
$ws.Range("J74").Value = 1464.9231
$ws.Range("K74").Value = 11906162
$ws.Range("L74").Value = 1464.9231
$ws.Range("M74").Value = -11905288
$ws.Range("N74").Value = -3212.9231

# row 77
$ws.Range("H77").Value = 9092324
$ws.Range("I77").Value = 11906162
$ws.Range("J77").Value = 1464.9231
$ws.Range("K77").Value = 59530810
$ws.Range("L77").Value = 7324.6155
$ws.Range("M77").Value = -59526442
$ws.Range("N77").Value = -16060.6155

# row 97
$ws.Range("H97").Value = 4839.5356
$ws.Range("I97").Value = 4655.8696
$ws.Range("J97").Value = 5684.4
$ws.Range("K97").Value = 4655.8696
$ws.Range("L97").Value = 5684.4
$ws.Range("M97").Value = -4159.8696
$ws.Range("N97").Value = -6676.4

# row 102
$ws.Range("H102").Value = 3046.3635
$ws.Range("I102").Value = 2851
$ws.Range("J102").Value = 5000
$ws.Range("K102").Value = 2851
$ws.Range("L102").Value = 5000
$ws.Range("M102").Value = -1229
$ws.Range("N102").Value = -8244

# row 132
$ws.Range("H132").Value = 9618223
$ws.Range("I132").Value = 13160544
$ws.Range("J132").Value = 3350.8572
$ws.Range("K132").Value = 39481632
$ws.Range("L132").Value = 10052.5716
$ws.Range("M132").Value = -39479102
$ws.Range("N132").Value = -15112.5716

# row 136
$ws.Range("H136").Value = 13159577
$ws.Range("I136").Value = 15153317
$ws.Range("J136").Value = 892.4
$ws.Range("K136").Value = 45459951
$ws.Range("L136").Value = 2677.2
$ws.Range("M136").Value = -45457401
$ws.Range("N136").Value = -7777.2

$ws = $wb.Worksheets.Item("BSM")
# row 107
$ws.Range("H107").Value = 1253.9231
$ws.Range("I107").Value = 1000.1
$ws.Range("J107").Value = 2100
$ws.Range("K107").Value = 1000.1
$ws.Range("L107").Value = 2100
$ws.Range("M107").Value = 919.9
$ws.Range("N107").Value = -5940

# row 134
$ws.Range("H134").Value = 2062.397
$ws.Range("I134").Value = 1301.9056
$ws.Range("J134").Value = 4749.467
$ws.Range("K134").Value = 3905.7168
$ws.Range("L134").Value = 14248.401
$ws.Range("M134").Value = -1370.7168
$ws.Range("N134").Value = -19318.401

$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 7577950.5
$ws.Range("I31").Value = 1987.5122
$ws.Range("J31").Value = 111116110
$ws.Range("K31").Value = 1987.5122
$ws.Range("L31").Value = 111116110
$ws.Range("M31").Value = -1692.5122
$ws.Range("N31").Value = -111116700

# row 34
$ws.Range("H34").Value = 7577950.5
$ws.Range("I34").Value = 1987.5122
$ws.Range("J34").Value = 111116110
$ws.Range("K34").Value = 1987.5122
$ws.Range("L34").Value = 111116110
$ws.Range("M34").Value = -1785.5122
$ws.Range("N34").Value = -111116514

# row 58
$ws.Range("H58").Value = 1496.7567
$ws.Range("I58").Value = 864.1177
$ws.Range("J58").Value = 8666.666999999999
$ws.Range("K58").Value = 864.1177
$ws.Range("L58").Value = 8666.666999999999
$ws.Range("M58").Value = -661.1177
$ws.Range("N58").Value = -9072.666999999999

# row 105
$ws.Range("H105").Value = 1880
$ws.Range("I105").Value = 1880
$ws.Range("K105").Value = 1880
$ws.Range("M105").Value = -133

# row 132
$ws.Range("H132").Value = 16130979
$ws.Range("I132").Value = 17242980
$ws.Range("K132").Value = 51728940
$ws.Range("M132").Value = -51726410

# row 134
$ws.Range("H134").Value = 1506.1923
$ws.Range("I134").Value = 1399.5209
$ws.Range("J134").Value = 2786.25
$ws.Range("K134").Value = 4198.5627
$ws.Range("L134").Value = 8358.75
$ws.Range("M134").Value = -1663.5627
$ws.Range("N134").Value = -13428.75

# row 136
$ws.Range("H136").Value = 1496.7567
$ws.Range("I136").Value = 864.1177
$ws.Range("J136").Value = 8666.666999999999
$ws.Range("K136").Value = 2592.3531
$ws.Range("L136").Value = 26000.001
$ws.Range("M136").Value = -42.35310000000027
$ws.Range("N136").Value = -31100.001

$ws = $wb.Worksheets.Item("CUL")
# row 55
$ws.Range("H55").Value = 520.4
$ws.Range("I55").Value = 125
$ws.Range("J55").Value = 595.7143
$ws.Range("K55").Value = 375
$ws.Range("L55").Value = 1787.1429
$ws.Range("M55").Value = -198
$ws.Range("N55").Value = -2141.1429

# row 107
$ws.Range("H107").Value = 1078.138
$ws.Range("J107").Value = 1273.5834
$ws.Range("L107").Value = 3820.7502
$ws.Range("N107").Value = -7660.7502

# row 113
$ws.Range("H113").Value = 682.4091
$ws.Range("I113").Value = 431.64102
$ws.Range("J113").Value = 882
$ws.Range("K113").Value = 1294.92306
$ws.Range("L113").Value = 2646
$ws.Range("M113").Value = 875.0769399999999
$ws.Range("N113").Value = -6986

# row 131
$ws.Range("H131").Value = 824.37897
$ws.Range("J131").Value = 859.4942600000001
$ws.Range("L131").Value = 2578.48278
$ws.Range("N131").Value = -12658.48278

$ws = $wb.Worksheets.Item("GSM")
# row 64
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

# row 67
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

# row 107
$ws.Range("H107").Value = 1679.2667
$ws.Range("I107").Value = 2725
$ws.Range("J107").Value = 484.14285
$ws.Range("K107").Value = 2725
$ws.Range("L107").Value = 484.14285
$ws.Range("M107").Value = -805
$ws.Range("N107").Value = -4324.14285

# row 126
$ws.Range("H126").Value = 3432.7407
$ws.Range("I126").Value = 1880.8334
$ws.Range("J126").Value = 4674.2666
$ws.Range("K126").Value = 5642.5002
$ws.Range("L126").Value = 14022.7998
$ws.Range("M126").Value = -3172.5002
$ws.Range("N126").Value = -18962.7998

# row 132
$ws.Range("H132").Value = 2953.1772
$ws.Range("I132").Value = 2159.709
$ws.Range("J132").Value = 4771.5415
$ws.Range("K132").Value = 6479.126999999999
$ws.Range("L132").Value = 14314.6245
$ws.Range("M132").Value = -3949.126999999999
$ws.Range("N132").Value = -19374.6245

$ws = $wb.Worksheets.Item("LTW")
# row 22
$ws.Range("H22").Value = 1264.4073
$ws.Range("I22").Value = 416.25
$ws.Range("J22").Value = 1621.5264
$ws.Range("K22").Value = 416.25
$ws.Range("L22").Value = 1621.5264
$ws.Range("M22").Value = -121.25
$ws.Range("N22").Value = -2211.5264

# row 27
$ws.Range("H27").Value = 1264.4073
$ws.Range("I27").Value = 416.25
$ws.Range("J27").Value = 1621.5264
$ws.Range("K27").Value = 416.25
$ws.Range("L27").Value = 1621.5264
$ws.Range("M27").Value = -309.25
$ws.Range("N27").Value = -1835.5264

# row 100
$ws.Range("H100").Value = 2071.12
$ws.Range("I100").Value = 2060
$ws.Range("J100").Value = 2078.5334
$ws.Range("K100").Value = 2060
$ws.Range("L100").Value = 2078.5334
$ws.Range("M100").Value = -1519
$ws.Range("N100").Value = -3160.5334

# row 122
$ws.Range("H122").Value = 4691.65
$ws.Range("I122").Value = 4461.7036
$ws.Range("K122").Value = 13385.1108
$ws.Range("M122").Value = -10935.1108

# row 132
$ws.Range("H132").Value = 7581993
$ws.Range("I132").Value = 3561.4119
$ws.Range("J132").Value = 33348660
$ws.Range("K132").Value = 10684.2357
$ws.Range("L132").Value = 100045980
$ws.Range("M132").Value = -8154.235700000001
$ws.Range("N132").Value = -100051040

# row 136
$ws.Range("H136").Value = 10003174
$ws.Range("I136").Value = 11364698
$ws.Range("J136").Value = 18667.5
$ws.Range("K136").Value = 34094094
$ws.Range("L136").Value = 56002.5
$ws.Range("M136").Value = -34091544
$ws.Range("N136").Value = -61102.5

$ws = $wb.Worksheets.Item("WVR")
# row 96
$ws.Range("H96").Value = 3350
$ws.Range("I96").Value = 2513.6365
$ws.Range("J96").Value = 4500
$ws.Range("K96").Value = 2513.6365
$ws.Range("L96").Value = 4500
$ws.Range("M96").Value = -1140.6365
$ws.Range("N96").Value = -7246

# row 126
$ws.Range("H126").Value = 3284.476
$ws.Range("I126").Value = 2421.4119
$ws.Range("K126").Value = 7264.2357
$ws.Range("M126").Value = -4794.2357

# row 132
$ws.Range("H132").Value = 1568.0944
$ws.Range("I132").Value = 1430.1025
$ws.Range("J132").Value = 1952.5
$ws.Range("K132").Value = 4290.3075
$ws.Range("L132").Value = 5857.5
$ws.Range("M132").Value = -1760.3075
$ws.Range("N132").Value = -10917.5

# row 136
$ws.Range("H136").Value = 915.525
$ws.Range("I136").Value = 678.0526
$ws.Range("J136").Value = 5427.5
$ws.Range("K136").Value = 2034.1578
$ws.Range("L136").Value = 16282.5
$ws.Range("M136").Value = 515.8422
$ws.Range("N136").Value = -21382.5
